$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price/volume figures (and a few coin rows that got reordered).
# Cell values that are plain text (names, URLs, percentage strings) can be
# assigned directly - Excel will not try to reinterpret them as numbers.
$textUpdates = @{
    "D2" = "30.058.90"
    "E2" = "  -0.20%  "
    "D3" = "1.912.34"
    "E3" = "  +0.11%  "
    "E4" = "  +0.03%  "
    "E5" = "  +8.55%  "
    "E6" = "  -1.00%  "
    "E7" = "  -0.02%  "
    "E8" = "  +3.83%  "
    "E9" = "  -0.62%  "
    "E10" = "  -0.44%  "
    "E11" = "  -0.55%  "
    "E12" = "  -2.67%  "
    "D13" = "1.900.50"
    "E13" = "  +0.22%  "
    "E14" = "  -2.33%  "
    "E15" = "  +1.57%  "
    "D16" = "30.064.85"
    "E16" = "  -0.18%  "
    "E17" = "  -2.48%  "
    "B18" = "Uniswap"
    "C18" = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
    "E18" = "  -0.07%  "
    "B19" = "BitcoinCash"
    "C19" = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
    "E19" = "  +3.51%  "
    "E20" = "  -0.62%  "
    "E21" = "  -0.01%  "
    "D22" = "2.149.43"
    "E22" = "  -1.65%  "
    "E23" = "  +0.06%  "
    "E24" = "  -1.82%  "
    "E25" = "  +0.89%  "
    "E26" = "  -1.22%  "
    "E27" = "  +9.48%  "
    "E28" = "  +0.14%  "
    "E29" = "  +0.07%  "
    "E30" = "  +3.08%  "
    "E31" = "  -2.00%  "
    "E32" = "  +0.31%  "
    "E33" = "  +0.74%  "
    "E34" = "  +3.33%  "
    "E35" = "  -3.93%  "
    "E36" = "  -1.47%  "
    "E37" = "  -0.04%  "
    "E38" = "  -1.28%  "
    "E39" = "  -0.32%  "
    "B40" = "TheSandbox"
    "C40" = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
    "E40" = "  -1.07%  "
    "B41" = "FraxShare"
    "C41" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
    "E41" = "  -3.03%  "
    "E42" = "  -2.63%  "
    "E43" = "  -0.08%  "
    "E44" = "  -4.29%  "
    "E45" = "  -1.22%  "
    "E46" = "  -2.14%  "
    "B47" = "Quant"
    "C47" = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
    "E47" = "  -1.15%  "
    "B48" = "EnergySwap"
    "C48" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "E48" = "  -1.43%  "
    "D49" = "2.056.05"
    "E49" = "  -0.73%  "
    "E50" = "  +3.73%  "
    "E51" = "  -1.13%  "
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# Price strings that look like plain numbers ("1.001", "26.37", ...) must stay
# as text (the source data always stored Price as a string). Force the cell to
# text format before writing, then clear the formatting again so the cell ends
# up with no explicit style, matching how the rest of the sheet is stored.
$numericLookingUpdates = @{
    "D5" = "0.8018"
    "D6" = "242.77"
    "D7" = "1.001"
    "D8" = "0.3214"
    "D9" = "26.37"
    "D10" = "0.06957"
    "D11" = "0.08023"
    "D12" = "0.7508"
    "D14" = "5.235"
    "D15" = "93.59"
    "D17" = "14.02"
    "D18" = "6.004"
    "D19" = "248.80"
    "D20" = "0.000007822"
    "D21" = "1.000"
    "D23" = "1.002"
    "D24" = "6.992"
    "D25" = "168.69"
    "D26" = "9.306"
    "D27" = "0.1409"
    "D28" = "19.00"
    "D29" = "2.057"
    "D30" = "1.394"
    "D31" = "1.524"
    "D32" = "4.358"
    "D33" = "4.129"
    "D34" = "0.05334"
    "D35" = "1.264"
    "D36" = "0.7407"
    "D37" = "2.733"
    "D38" = "0.01928"
    "D39" = "2.791"
    "D40" = "0.4468"
    "D41" = "6.155"
    "D42" = "72.74"
    "D43" = "1.001"
    "D44" = "1.906"
    "D45" = "0.8300"
    "D46" = "7.615"
    "D47" = "100.72"
    "D48" = "9.823"
    "D50" = "963.24"
    "D51" = "36.50"
}
foreach ($addr in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$addr]
    $cell.ClearFormats()
}

Write-Host "Applied cryptos update"
